$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Standard CSRp")

# Relabel the "Time (ms)" sub-header to "Time (us)" since the new rows below
# are reported in microseconds.
$ws.Range("G5").Value = "Time (us)"

# Add a few more data rows (8, 16 and 32 pipes) to the speed-test columns.
$ws.Range("G9").Value = 0.019287
$ws.Range("H9").Value = 0.103695
$ws.Range("I9").Formula = "=H9/H6"

$ws.Range("G10").Value = 0.010292
$ws.Range("H10").Value = 0.194331
$ws.Range("I10").Formula = "=H10/H6"

$ws.Range("G11").Value = 0.005778
$ws.Range("H11").Value = 0.346162
$ws.Range("I11").Formula = "=H11/H6"

# Keep the selection/zoom in sync with where the author was last looking.
$ws.Range("I9").Select()
$excel.ActiveWindow.Zoom = 89
